# Updates the cryptos list: Price (column D) and Volume(1h) (column E)
# values for the rows that changed, per the upstream data refresh.
#
# Values in columns D and E are stored as text (not numbers) in the
# workbook so that things like "522.00", "1.981.02" or "  +3.52%  "
# round-trip byte-for-byte. Excel's COM layer auto-coerces numeric-
# looking strings assigned to Range.Value into real numbers, so we
# temporarily force the cell's number format to Text ("@") before the
# assignment and restore the original style afterwards (leaving the
# cell's visual style/format exactly as it was before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue $ws 'D2' '59.355.32'
Set-TextValue $ws 'E2' '  +3.68%  '
Set-TextValue $ws 'D3' '2.594.86'
Set-TextValue $ws 'E3' '  +2.21%  '
Set-TextValue $ws 'E4' '  +0.10%  '
Set-TextValue $ws 'D5' '522.00'
Set-TextValue $ws 'E5' '  +1.45%  '
Set-TextValue $ws 'D6' '140.58'
Set-TextValue $ws 'E6' '  +0.55%  '
Set-TextValue $ws 'E7' '  -0.34%  '
Set-TextValue $ws 'E8' '  +1.84%  '
Set-TextValue $ws 'D9' '2.617.59'
Set-TextValue $ws 'E9' '  +3.12%  '
Set-TextValue $ws 'E10' '  +0.80%  '
Set-TextValue $ws 'E11' '  +2.05%  '
Set-TextValue $ws 'D12' '0.332'
Set-TextValue $ws 'E12' '  +2.45%  '
Set-TextValue $ws 'E13' '  +1.95%  '
Set-TextValue $ws 'D14' '3.057.18'
Set-TextValue $ws 'E14' '  +2.48%  '
Set-TextValue $ws 'D15' '59.317.60'
Set-TextValue $ws 'E15' '  +3.64%  '
Set-TextValue $ws 'D16' '20.40'
Set-TextValue $ws 'E16' '  +2.08%  '
Set-TextValue $ws 'D17' '2.608.08'
Set-TextValue $ws 'E17' '  +3.07%  '
Set-TextValue $ws 'D18' '0.0000133'
Set-TextValue $ws 'E18' '  +0.18%  '
Set-TextValue $ws 'D19' '338.57'
Set-TextValue $ws 'E19' '  +2.11%  '
Set-TextValue $ws 'D20' '4.32'
Set-TextValue $ws 'E20' '  +1.29%  '
Set-TextValue $ws 'D21' '10.19'
Set-TextValue $ws 'E21' '  +1.14%  '
Set-TextValue $ws 'D22' '6.50'
Set-TextValue $ws 'E22' '  +6.30%  '
Set-TextValue $ws 'E23' '  -0.28%  '
Set-TextValue $ws 'D24' '66.41'
Set-TextValue $ws 'E24' '  +3.22%  '
Set-TextValue $ws 'D25' '0.168'
Set-TextValue $ws 'E25' '  +1.30%  '
Set-TextValue $ws 'D26' '0.404'
Set-TextValue $ws 'E26' '  +0.98%  '
Set-TextValue $ws 'D27' '0.995'
Set-TextValue $ws 'E27' '  -0.34%  '
Set-TextValue $ws 'D28' '7.06'
Set-TextValue $ws 'E28' '  +1.84%  '
Set-TextValue $ws 'E29' '  -0.15%  '
Set-TextValue $ws 'E30' '  -3.05%  '
Set-TextValue $ws 'E31' '  -4.88%  '
Set-TextValue $ws 'E32' '  +2.21%  '
Set-TextValue $ws 'E33' '  +0.94%  '
Set-TextValue $ws 'D34' '149.27'
Set-TextValue $ws 'E34' '  +0.34%  '
Set-TextValue $ws 'E35' '  +0.93%  '
Set-TextValue $ws 'E36' '  +0.13%  '
Set-TextValue $ws 'D37' '36.37'
Set-TextValue $ws 'E37' '  +1.84%  '
Set-TextValue $ws 'D38' '1.46'
Set-TextValue $ws 'E38' '  +3.88%  '
Set-TextValue $ws 'D39' '0.833'
Set-TextValue $ws 'E39' '  +1.26%  '
Set-TextValue $ws 'E40' '  -2.43%  '
Set-TextValue $ws 'E41' '  +2.17%  '
Set-TextValue $ws 'D42' '0.995'
Set-TextValue $ws 'E42' '  -0.48%  '
Set-TextValue $ws 'D43' '275.88'
Set-TextValue $ws 'E43' '  +6.87%  '
Set-TextValue $ws 'D44' '10.74'
Set-TextValue $ws 'E44' '  +1.31%  '
Set-TextValue $ws 'D46' '0.0954'
Set-TextValue $ws 'E46' '  +0.30%  '
Set-TextValue $ws 'E47' '  +0.01%  '
Set-TextValue $ws 'D48' '18.58'
Set-TextValue $ws 'E48' '  +1.26%  '
Set-TextValue $ws 'D49' '1.980.62'
Set-TextValue $ws 'E49' '  +0.72%  '
Set-TextValue $ws 'E50' '  +2.80%  '
Set-TextValue $ws 'D51' '0.0221'
Set-TextValue $ws 'E51' '  -0.10%  '
